# Combine the three runs that render the "input date" header cell's
# conditional template text into a single run whose text applies a
# "date" filter (with a Swiss dd.MM.YYYY format) to both the electronic
# and paper input-date placeholders.
#
# Old (3 runs):
#   {% if inputDateHeader %}{{ inputDateHeader }}
#   {% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}
#   {% else %}-{% endif %}
#
# New (1 run):
#   {% if inputDateHeader %}{{ inputDateHeader | date("dd.MM.YYYY") }}{% if paperInputDateHeader %} ({{ paperInputDateHeader | date("dd.MM.YYYY") }}){% else %}{% endif %}{% else %}-{% endif %}

$d = $word.ActiveDocument

$old = "{% if inputDateHeader %}{{ inputDateHeader }}{% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}{% else %}-{% endif %}"
$new = '{% if inputDateHeader %}{{ inputDateHeader | date("dd.MM.YYYY") }}{% if paperInputDateHeader %} ({{ paperInputDateHeader | date("dd.MM.YYYY") }}){% else %}{% endif %}{% else %}-{% endif %}'

# Locate the exact span of text across the three runs (Find matches
# across run boundaries on the Range's plain text).
$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the inputDateHeader template text to update"
}

# Assign straight to Range.Text (rather than Find.Execute's Replace
# argument) so Word doesn't run its "smart quotes" AutoFormat over the
# new `date("dd.MM.YYYY")` literal. This also collapses the matched
# span -- which originally spanned three separate <w:r> runs -- down to
# a single run that inherits the formatting of the first run.
$rng.Text = $new
